# chore: simulator full-month coverage, persist logs, fix employees
#
# 1) Fix mis-assigned client names on the "Weekly Timesheet" sheet (and the
#    mirrored rows on "Jason Schema").
# 2) Populate the simulator's Rate/Total columns (previously all zero) with
#    full-month coverage numbers, and roll the new totals up into the
#    subtotal / grand-total rows.
# 3) Re-issue Doug Kinsey's employee id.

$wb = $excel.ActiveWorkbook

$timesheet = $wb.Worksheets.Item("Weekly Timesheet")
$schema = $wb.Worksheets.Item("Jason Schema")

# --- 1) Client name corrections -------------------------------------------------

$timesheet.Range("B2").Value = "Hunter"
$timesheet.Range("B3").Value = "Tubergen"
$timesheet.Range("B4").Value = "Field"
$timesheet.Range("B5").Value = "Bottomley"
$timesheet.Range("B6").Value = "Zygmunt"

$schema.Range("D2").Value = "Hunter"
$schema.Range("D3").Value = "Tubergen"
$schema.Range("D4").Value = "Field"
$schema.Range("D5").Value = "Bottomley"
$schema.Range("D6").Value = "Zygmunt"

# --- 2) Rate / Total figures ----------------------------------------------------

$rows = 2, 3, 4, 5, 6
foreach ($r in $rows) {
    $timesheet.Range("E$r").Value = 92
    $timesheet.Range("F$r").Value = 736

    $schema.Range("F$r").Value = 92
    $schema.Range("G$r").Value = 736
}

$timesheet.Range("F8").Value = 3680
$timesheet.Range("F11").Value = 3680
$timesheet.Range("F13").Value = 3680

# --- 3) Employee id rotation -----------------------------------------------------

foreach ($r in $rows) {
    $schema.Range("B$r").Value = "emp_yde33znx"
}
